$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column P: header year 2022 (row 3) ---
$ws.Range("P3").Value = 2022

# Copy the formatting (font/border/fill/alignment) of the existing
# neighbouring column O cells onto the new column P cells for the rows
# that don't carry the new "0.0" numeric format (the thin/thick-border
# spacer row and the year-header row). Row 1 is left untouched -
# column P never gets a cell there, same as the source workbook.
$ws.Range("O2").Copy() | Out-Null
$ws.Range("P2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("O3").Copy() | Out-Null
$ws.Range("P3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- New column P data rows 4-6 ---
$ws.Range("P5").Value = 1339.6
$ws.Range("P6").Value = 6300.5
$ws.Range("P4").Formula = "=P5/P6*1000"

# Give the three new data cells the same formatting as their row
# neighbours in column O before the number-format change below is
# applied uniformly across the whole D:P data block.
$ws.Range("O4").Copy() | Out-Null
$ws.Range("P4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("O5").Copy() | Out-Null
$ws.Range("P5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("O6").Copy() | Out-Null
$ws.Range("P6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Re-format the whole numeric data block (now including column P)
# from two decimals ("0.00") to one decimal ("0.0"). This both
# registers the new custom number format and normalises O5 (which had
# drifted to a plain/general format) back in line with the rest of its
# row, exactly like the rest of the row/column.
$ws.Range("D4:L4").NumberFormat = "0.0"
$ws.Range("M4:P4").NumberFormat = "0.0"
$ws.Range("D5:P5").NumberFormat = "0.0"
$ws.Range("D6:P6").NumberFormat = "0.0"

# --- Selection moves off the old P16 anchor ---
$ws.Range("S4").Select() | Out-Null
